# Season up to 1/17
$wb = $excel.ActiveWorkbook

$gamesWs = $wb.Worksheets.Item("Games")
$nextWs = $wb.Worksheets.Item("Next")

# 1. Append the now-played game (2024-01-15, DAL) to the "Games" sheet as row 42,
#    using the box-score values from the result.
$newRow = 42
$gamesWs.Cells.Item($newRow, 1).Value = 41
$gamesWs.Cells.Item($newRow, 2).Value = 45306
$gamesWs.Cells.Item($newRow, 2).NumberFormat = "YYYY-MM-DD"
$gamesWs.Cells.Item($newRow, 3).Value = -1
$gamesWs.Cells.Item($newRow, 4).Value = 120
$gamesWs.Cells.Item($newRow, 5).Value = 96
$gamesWs.Cells.Item($newRow, 6).Value = 0.549
$gamesWs.Cells.Item($newRow, 7).Value = 8.300000000000001
$gamesWs.Cells.Item($newRow, 8).Value = 21.6
$gamesWs.Cells.Item($newRow, 9).Value = 0.383
$gamesWs.Cells.Item($newRow, 10).Value = 125.1
$gamesWs.Cells.Item($newRow, 11).Value = "DAL"
$gamesWs.Cells.Item($newRow, 12).Value = 125
$gamesWs.Cells.Item($newRow, 13).Value = 0.5649999999999999
$gamesWs.Cells.Item($newRow, 14).Value = 9.9
$gamesWs.Cells.Item($newRow, 15).Value = 31
$gamesWs.Cells.Item($newRow, 16).Value = 0.341
$gamesWs.Cells.Item($newRow, 17).Value = 130.3
$gamesWs.Cells.Item($newRow, 18).Value = 0
$gamesWs.Cells.Item($newRow, 19).Value = 0

# 2. Remove the now-played game from the "Next" sheet (row 2), shifting the
#    remaining upcoming games up by one row.
$nextWs.Rows.Item(2).Delete()
